$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fall 23 week 6 inputs - update matchup average values
$ws.Range("F3").Value = 1.21
$ws.Range("E4").Value = 1.23
$ws.Range("C5").Value = 1.38
$ws.Range("D5").Value = 1.35
$ws.Range("G5").Value = 0.76
$ws.Range("C6").Value = 1.47
$ws.Range("D7").Value = 1.78
$ws.Range("E7").Value = 1.89
